$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 135 ("「世界の七不思議」" post) - this shifts every
# row below it up by one, matching the renumbering seen in the diff.
$ws.Rows.Item(135).Delete()
